# Applies the cryptos.xlsx data refresh described by the commit
# "Updated cryptos list on Thu Sep 12 21:40:04 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell reference -> new text value
$updates = @{
    'D2' = '58.051.07'
    'E2' = '  +2.04%  '
    'D3' = '2.347.98'
    'E3' = '  +0.28%  '
    'E4' = '  -0.10%  '
    'D5' = '539.83'
    'E5' = '  +1.95%  '
    'D6' = '135.81'
    'E6' = '  +2.64%  '
    'D7' = '0.999'
    'E7' = '  +0.24%  '
    'D8' = '0.565'
    'E8' = '  +5.69%  '
    'E9' = '  +0.48%  '
    'D10' = '5.56'
    'E10' = '  +4.92%  '
    'E11' = '  -0.69%  '
    'E12' = '  +1.79%  '
    'D13' = '23.80'
    'E13' = '  +1.56%  '
    'D14' = '2.763.07'
    'E14' = '  +1.03%  '
    'D15' = '57.971.13'
    'E15' = '  +1.86%  '
    'E16' = '  +0.71%  '
    'D17' = '2.337.16'
    'E17' = '  +0.10%  '
    'D18' = '10.70'
    'E18' = '  +2.69%  '
    'D19' = '331.81'
    'E19' = '  -1.30%  '
    'D20' = '4.28'
    'E20' = '  +2.91%  '
    'D21' = '6.79'
    'E21' = '  -0.97%  '
    'D22' = '0.998'
    'E22' = '  -0.08%  '
    'D23' = '62.87'
    'E23' = '  +2.14%  '
    'E24' = '  -0.29%  '
    'D25' = '8.50'
    'E25' = '  -2.26%  '
    'D26' = '0.998'
    'E26' = '  +0.45%  '
    'D27' = '1.38'
    'E27' = '  +2.29%  '
    'E28' = '  +1.74%  '
    'D29' = '171.92'
    'E29' = '  -0.36%  '
    'D30' = '0.0₃0737'
    'E30' = '  +1.70%  '
    'D31' = '6.14'
    'E31' = '  +0.68%  '
    'D32' = '1.03'
    'E32' = '  +11.20%  '
    'D33' = '18.44'
    'E33' = '  -0.20%  '
    'B35' = 'NEARProtocol'
    'C35' = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
    'D35' = '4.23'
    'E35' = '  +6.23%  '
    'B36' = 'FirstDigitalUSD'
    'C36' = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
    'D36' = '1.00'
    'E36' = '  +0.57%  '
    'E37' = '  -0.20%  '
    'E38' = '  +4.73%  '
    'D39' = '39.21'
    'E39' = '  +0.34%  '
    'D40' = '145.27'
    'E40' = '  -2.69%  '
    'D41' = '293.52'
    'E41' = '  +4.06%  '
    'D42' = '0.378'
    'E42' = '  +0.72%  '
    'D43' = '3.65'
    'E43' = '  +1.17%  '
    'D44' = '0.0948'
    'E44' = '  +1.94%  '
    'D45' = '19.33'
    'E45' = '  +2.51%  '
    'D46' = '0.0504'
    'E46' = '  +0.75%  '
    'D47' = '0.563'
    'E47' = '  +0.85%  '
    'E48' = '  +1.42%  '
    'B49' = 'EnergySwap'
    'C49' = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
    'D49' = '17.51'
    'E49' = '  +0.32%  '
    'B50' = 'Polygon'
    'C50' = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
    'D50' = '0.381'
    'E50' = '  -0.11%  '
    'D51' = '11.05'
    'E51' = '  +0.38%  '
}

# Cells whose new value looks like a plain number (e.g. "539.83").
# These must be forced to stay text so Excel does not coerce them into
# numeric cells (which would also lose formatting like trailing zeros).
$forceText = @(
    'D5', 'D6', 'D7', 'D8', 'D10', 'D13', 'D18', 'D19', 'D20', 'D21',
    'D22', 'D23', 'D25', 'D26', 'D27', 'D29', 'D31', 'D32', 'D33', 'D35',
    'D36', 'D39', 'D40', 'D41', 'D42', 'D43', 'D44', 'D45', 'D46', 'D47',
    'D49', 'D50', 'D51'
)

foreach ($ref in $updates.Keys) {
    $rng = $ws.Range($ref)
    if ($forceText -contains $ref) {
        $origStyle = $rng.Style
        $rng.NumberFormat = "@"
        $rng.Value = $updates[$ref]
        $rng.Style = $origStyle
    } else {
        $rng.Value = $updates[$ref]
    }
}
